# This run represents a fresh scrape cycle of KHL injury data:
#  - every row on the "snapshot" sheet gets a refreshed scraped_at timestamp
#  - Savitsky Kirill (Barys) is newly injured -> inserted into "snapshot"
#    and logged on "new_injured"
#  - Kalinichenko Roman (Sibir) has recovered -> removed from "snapshot"
#    and logged on "returned"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "snapshot" sheet: insert the new injury row, drop the returned
#    player's row, then refresh every scraped_at (column K) value.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("snapshot")

# Insert a fresh blank row for Savitsky Kirill right after Boyarkin Nikita
# (row 9), pushing Waterspoon Tyler and everyone below down by one.
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = "БАР"
$ws.Cells.Item(10, 2).Value = "Барыс"
$ws.Cells.Item(10, 3).Value = "barys"
$ws.Cells.Item(10, 4).Value = "Савицкий Кирилл"
$ws.Cells.Item(10, 5).Value = "'84"
$ws.Cells.Item(10, 6).Value = "нападающий"
$ws.Cells.Item(10, 7).Value = "'17901"
$ws.Cells.Item(10, 8).Value = "1369_БАР_савицкийкирилл"
$ws.Cells.Item(10, 9).Value = "injured_active"
$ws.Cells.Item(10, 10).Value = "https://www.khl.ru/clubs/barys/team/"
$ws.Cells.Item(10, 11).Value = "2025-11-19T10:01:28.356049+00:00"

# After the insert above, Kalinichenko Roman's row (originally row 19)
# now sits at row 20 - remove it, he's back from injury.
$ws.Rows.Item(20).Delete()

# Refresh the scraped_at timestamp (column K) for every remaining data
# row (2-40) to reflect this scrape pass.
$scrapedAt = @(
  "2025-11-19T10:01:07.669188+00:00",
  "2025-11-19T10:01:07.669233+00:00",
  "2025-11-19T10:01:07.669316+00:00",
  "2025-11-19T10:01:12.752236+00:00",
  "2025-11-19T10:01:12.752269+00:00",
  "2025-11-19T10:01:17.774573+00:00",
  "2025-11-19T10:01:23.267471+00:00",
  "2025-11-19T10:01:28.356018+00:00",
  "2025-11-19T10:01:28.356049+00:00",
  "2025-11-19T10:01:28.356068+00:00",
  "2025-11-19T10:01:38.456213+00:00",
  "2025-11-19T10:01:43.553855+00:00",
  "2025-11-19T10:01:48.568090+00:00",
  "2025-11-19T10:01:48.568125+00:00",
  "2025-11-19T10:01:48.568146+00:00",
  "2025-11-19T10:01:53.536030+00:00",
  "2025-11-19T10:01:58.571368+00:00",
  "2025-11-19T10:02:03.602987+00:00",
  "2025-11-19T10:02:08.617446+00:00",
  "2025-11-19T10:02:08.617480+00:00",
  "2025-11-19T10:02:08.617500+00:00",
  "2025-11-19T10:02:13.763035+00:00",
  "2025-11-19T10:02:13.763070+00:00",
  "2025-11-19T10:02:18.797388+00:00",
  "2025-11-19T10:02:18.797423+00:00",
  "2025-11-19T10:02:18.797447+00:00",
  "2025-11-19T10:02:24.223155+00:00",
  "2025-11-19T10:02:24.223191+00:00",
  "2025-11-19T10:02:29.257514+00:00",
  "2025-11-19T10:02:29.257547+00:00",
  "2025-11-19T10:02:29.257570+00:00",
  "2025-11-19T10:02:29.257591+00:00",
  "2025-11-19T10:02:29.257609+00:00",
  "2025-11-19T10:02:34.771679+00:00",
  "2025-11-19T10:02:34.771717+00:00",
  "2025-11-19T10:02:44.866643+00:00",
  "2025-11-19T10:02:44.866676+00:00",
  "2025-11-19T10:02:44.866694+00:00",
  "2025-11-19T10:02:49.904008+00:00"
)
for ($i = 0; $i -lt $scrapedAt.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $scrapedAt[$i]
}

# ---------------------------------------------------------------------
# 2) "returned" sheet: log Kalinichenko Roman coming back from injury.
# ---------------------------------------------------------------------
$returned = $wb.Worksheets.Item("returned")
$returned.Cells.Item(3, 1).Value = "СИБ"
$returned.Cells.Item(3, 2).Value = "Сибирь"
$returned.Cells.Item(3, 3).Value = "Калиниченко Роман"
$returned.Cells.Item(3, 4).Value = "1369_СИБ_калиниченкороман"
$returned.Cells.Item(3, 5).Value = "RETURN"
$returned.Cells.Item(3, 6).Value = "2025-11-19T18:02:50.406848+08:00"
$returned.Cells.Item(3, 7).Value = "'2025-11-19"

# ---------------------------------------------------------------------
# 3) "new_injured" sheet: log Savitsky Kirill's new injury.
# ---------------------------------------------------------------------
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Cells.Item(2, 1).Value = "БАР"
$newInjured.Cells.Item(2, 2).Value = "Барыс"
$newInjured.Cells.Item(2, 3).Value = "Савицкий Кирилл"
$newInjured.Cells.Item(2, 4).Value = "1369_БАР_савицкийкирилл"
$newInjured.Cells.Item(2, 5).Value = "INJURED_NEW"
$newInjured.Cells.Item(2, 6).Value = "2025-11-19T18:02:50.406848+08:00"
$newInjured.Cells.Item(2, 7).Value = "'2025-11-19"
